# Update the Rally credentials row (row 2) on Sheet1:
#  - W2 (RALLY_EMAIL): new login email, now shown/styled as a hyperlink
#  - X2 (RALLY_PASSWORD): new login password
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("W2").Value = "abhinay.marapaka@rallyhealth.com"
$ws.Range("W2").Style = "Hyperlink"

$ws.Range("X2").Value = "AbhinayElias22$"

[void]$ws.Range("X2").Select()
